$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "last row" date-only format currently on row 42 (A42) before
# we change it, so the new last row (43) can reuse it.
$lastRowDateFormat = $ws.Cells.Item(42, 1).NumberFormat

# Previously the last row (42) had the "last row" date-only format.
# Now that a new row is appended, row 42 becomes a regular data row, so its
# date format should match the rest of the column (same as A41).
$ws.Cells.Item(42, 1).NumberFormat = $ws.Cells.Item(41, 1).NumberFormat

# Append the new daily data row (row 43).
$ws.Cells.Item(43, 1).Value = 45783
$ws.Cells.Item(43, 2).Value = 179
$ws.Cells.Item(43, 3).Value = 182
$ws.Cells.Item(43, 4).Value = 177

# The new last row takes on the "last row" date-only format, matching what
# row 42 used to have.
$ws.Cells.Item(43, 1).NumberFormat = $lastRowDateFormat
